$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (no explicit) style, used to keep D-column cells
# unstyled after forcing them to be interpreted as text via the leading apostrophe.
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "'39.685.06"
$ws.Range("D2").Style = $plainStyle
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "'2.216.04"
$ws.Range("D3").Style = $plainStyle
$ws.Range("E3").Value = "  -5.29%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = $plainStyle
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'298.94"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  -3.74%  "
$ws.Range("D6").Value = "'83.78"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("D7").Value = "'0.515"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  -2.74%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.466"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "  -3.73%  "
$ws.Range("D10").Value = "'0.0781"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  -3.71%  "
$ws.Range("D11").Value = "'29.65"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").Value = "'46.23"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  -11.81%  "
$ws.Range("D13").Value = "'0.107"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").Value = "'2.557.88"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  -5.24%  "
$ws.Range("D15").Value = "'6.30"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  -2.17%  "
$ws.Range("D16").Value = "'14.13"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  -4.26%  "
$ws.Range("D17").Value = "'2.212.37"
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = "  -6.61%  "
$ws.Range("D18").Value = "'0.718"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "  -5.20%  "
$ws.Range("D19").Value = "'39.588.35"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").Value = "'0.0₃0879"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("D21").Value = "'5.74"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  -6.19%  "
$ws.Range("D22").Value = "'65.04"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  -4.29%  "
$ws.Range("D23").Value = "'10.42"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D24").Value = "'231.91"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "'2.43"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  -5.09%  "
$ws.Range("D27").Value = "'1.84"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("D28").Value = "'22.75"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("D29").Value = "'2.18"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  +2.51%  "
$ws.Range("D30").Value = "'9.17"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("D31").Value = "'32.33"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  -7.33%  "
$ws.Range("D32").Value = "'149.33"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "'4.84"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = "  -5.35%  "
$ws.Range("D35").Value = "'2.38"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("D36").Value = "'0.0702"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  -2.28%  "
$ws.Range("D37").Value = "'16.13"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  +3.07%  "
$ws.Range("D38").Value = "'0.111"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  -2.52%  "
$ws.Range("D39").Value = "'0.0971"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").Value = "'2.66"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  -5.37%  "
$ws.Range("D41").Value = "'1.65"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  -4.33%  "
$ws.Range("D42").Value = "'3.67"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  -5.60%  "
$ws.Range("D43").Value = "'1.927.66"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").Value = "'2.14"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  -3.45%  "
$ws.Range("D45").Value = "'0.0266"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").Value = "'9.23"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  -1.92%  "
$ws.Range("D47").Value = "'16.39"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  -7.75%  "
$ws.Range("D48").Value = "'2.61"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  -3.44%  "
$ws.Range("D49").Value = "'2.429.55"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  -5.00%  "
$ws.Range("D50").Value = "'70.92"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "'88.82"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  -4.28%  "

Write-Host "Updated cryptos list"